$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '58.037.69'
$ws.Range("E2").Value = '  +1.81%  '
$ws.Range("D3").Value = '3.133.55'
$ws.Range("E3").Value = '  +1.70%  '
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '534.93'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +3.04%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '138.98'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +2.56%  '
$ws.Range("E7").Value = '  -0.05%  '
$ws.Range("E8").Value = '  +10.89%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '7.35'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.18%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.109'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +2.39%  '
$ws.Range("E11").Value = '  +4.59%  '
$ws.Range("E12").Value = '  +2.95%  '
$ws.Range("D13").Value = '3.667.14'
$ws.Range("E13").Value = '  +1.53%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '25.73'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +2.00%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.0000169'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +5.13%  '
$ws.Range("D16").Value = '58.104.25'
$ws.Range("E16").Value = '  +1.74%  '
$ws.Range("D18").Value = '3.140.61'
$ws.Range("E18").Value = '  +1.80%  '
$ws.Range("E19").Value = '  +4.02%  '
$ws.Range("E20").Value = '  +3.95%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '374.68'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +7.94%  '
$ws.Range("E22").Value = '  +0.27%  '
$ws.Range("E23").Value = '  -0.66%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '70.12'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +2.73%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.512'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +3.03%  '
$ws.Range("E26").Value = '  +1.26%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.01'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.87%  '
$ws.Range("D28").Value = '0.0₃0890'
$ws.Range("E28").Value = '  +3.48%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.71'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +6.18%  '
$ws.Range("E30").Value = '  +6.71%  '
$ws.Range("E31").Value = '  +1.10%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '21.74'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +4.57%  '
$ws.Range("E33").Value = '  +5.96%  '
$ws.Range("E34").Value = '  +4.37%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '160.78'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.92%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '6.23'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +3.88%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.33'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +8.06%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '25.55'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.06%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.69'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +8.07%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0676'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +3.45%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '4.21'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +4.19%  '
$ws.Range("D42").Value = '2.558.92'
$ws.Range("E42").Value = '  +7.53%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '38.71'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +5.83%  '
$ws.Range("E44").Value = '  +1.23%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0272'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +3.10%  '
$ws.Range("E46").Value = '  -0.08%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '6.21'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +4.66%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.980'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +2.99%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0986'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +10.33%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '20.17'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +2.88%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.749'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.13%  '
